$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leaderboard")

$ws.Range("B2").Value = 700
$ws.Range("C2").Value = 560
$ws.Range("D2").Value = 6
$ws.Range("F2").Value = 6

$ws.Range("B3").Value = 256
$ws.Range("C3").Value = 162
$ws.Range("D3").Value = 75
$ws.Range("E3").Value = 14
